$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 217, shifting existing rows 217-250 down to 218-251.
$ws.Rows.Item(217).Insert()

# Populate the newly inserted row 217 with the new weekly record.
$ws.Range("A217").Value = 11
$ws.Range("B217").Value = "Vega Monumental Concepción"
$ws.Range("C217").Value = "Bíobío"
$ws.Range("D217").Value = 45127
$ws.Range("E217").Value = 8
$ws.Range("F217").Value = 100112032
$ws.Range("G217").Value = "Zapallo italiano"
$ws.Range("H217").Value = "Sin especificar"
$ws.Range("I217").Value = "Primera"
$ws.Range("J217").Value = 80
$ws.Range("K217").Value = 17000
$ws.Range("L217").Value = 17000
$ws.Range("M217").Value = 17000
$ws.Range("N217").Value = "$/caja 50 unidades"
$ws.Range("O217").Value = "Región de Arica y Parinacota"
$ws.Range("P217").Value = 340
$ws.Range("Q217").Value = 50
$ws.Range("R217").Value = "Hortaliza"
